# Fix typo "put put" -> "put" on slide 7, shape "parameters..."
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(5)
$tr = $shape.TextFrame.TextRange
$para2 = $tr.Paragraphs(2)
$run1 = $para2.Runs(1)
$run1.Text = "the placeholder(s) we put between the method’s parentheses when we "
